# Updates cryptos list values (price + volume) to match the latest scrape.
# Generated from the authoritative cell-by-cell diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.319.91"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.868.63"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'692.13"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "'173.31"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D7").Value = "3.864.59"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  +4.72%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +6.81%  "
$ws.Range("D14").Value = "'36.68"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "4.524.56"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "3.883.98"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "71.342.23"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'17.85"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").Value = "'496.76"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").Value = "'0.725"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'84.97"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("D26").Value = "'12.40"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").Value = "'10.61"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").Value = "4.030.29"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +10.53%  "
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").Value = "'29.85"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").Value = "'9.34"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "3.823.73"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.104"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("E40").Value = "  +13.80%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Value = "'6.07"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'1.02"
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D46").Value = "'164.40"
$ws.Range("E46").Value = "  +3.28%  "
$ws.Range("D47").Value = "'0.000307"
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'48.69"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'44.58"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = "  +2.59%  "
